$d = $word.ActiveDocument

# 1) Fix the split hyperlink text: "portfolio-" + "w" + "ebsite/" -> "portfolio-website/"
#    (Match starting one character inside the first run so the merged run keeps
#    its original rPr/rStyle=Hyperlink formatting instead of losing it.)
$d.Content.Find.Execute("ttps://limitdak.github.io/portfolio-website/", $true, $false, $false, $false, $false, $true, 1, $false, "ttps://limitdak.github.io/portfolio-website/", 2)

# 2) Fix the certification bullet text
$d.Content.Find.Execute("Finished all lessons within Responsive Web Design; nearing Certification", $true, $false, $false, $false, $false, $true, 1, $false, "Certification in Responsive Web Design", 2)
